# Add a new customer row (phone 79174418) with 0 total_points, mirroring
# the existing rows where the phone number is stored as text and the
# points start at 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phone numbers in column A are text, even though they look numeric
# (e.g. "79174418"), so format the cell as Text before typing the value
# -- otherwise Excel would auto-convert it to a number.
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "79174418"
# Drop back to the default cell style (the Text number-format was only
# needed to force text entry, not to keep a lasting custom format).
$ws.Range("A13").Style = "Normal"

# No birthday on file for this customer, so B13 stays blank.

# Starting point total for a brand-new customer.
$ws.Range("C13").Value = 0
